# Append the new 2025-10-03 profit row (A47:B47) to the sheet, matching the
# existing table layout (date stored as plain text, profit as a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column holds literal text like "10/02/2025", not real dates.
# Setting NumberFormat to Text ("@") before writing the value stops Excel's
# automatic date-recognition from turning the string into a date serial.
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "10/03/2025"
$ws.Range("B47").Value = 15474.33

# Re-apply the formatting of the preceding data row (A46) so the new cell
# ends up with the same (default) style as the rest of the column instead
# of keeping the temporary Text number format / quote-prefix flag.
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
